# Round the numeric data in B2:E13 to integer values (as whole numbers),
# matching the "write Ontpl_/Pot_ files to disk as integer data" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:E13")

foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = [Math]::Round([double]$val)
    }
}
